# Applies the "pridani nekterych komponent pro regulaci" change:
# - adds a new component row (display i2c) with its price/link
# - fixes VAT multiplier used in column E from 1.19 to 1.21 for rows 33-39
# - adds "Cena bez DPH celkem" (F) and "Cena s DPH celkem" (G) formulas for rows 33-39
# - updates the frozen-pane view / selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 33 ---
$ws.Range("E33").Formula = "=(D33*1.21)"
$ws.Range("F33").Formula = "=(D33*C33)"
$ws.Range("G33").Formula = "=(F33*1.21)"

# --- Row 34 ---
$ws.Range("E34").Formula = "=(D34*1.21)"
$ws.Range("F34").Formula = "=(D34*C34)"
$ws.Range("G34").Formula = "=(F34*1.21)"

# --- Row 35 ---
$ws.Range("E35").Formula = "=(D35*1.21)"
$ws.Range("F35").Formula = "=(D35*C35)"
$ws.Range("G35").Formula = "=(F35*1.21)"

# --- Row 36 ---
$ws.Range("E36").Formula = "=(D36*1.21)"
$ws.Range("F36").Formula = "=(D36*C36)"
$ws.Range("G36").Formula = "=(F36*1.21)"

# --- Row 37 ---
$ws.Range("E37").Formula = "=(D37*1.21)"
$ws.Range("F37").Formula = "=(D37*C37)"
$ws.Range("G37").Formula = "=(F37*1.21)"

# --- Row 38: new component "display i2c" ---
$ws.Range("B38").Value = "display i2c"
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 276
$ws.Range("E38").Formula = "=(D38*1.21)"
$ws.Range("F38").Formula = "=(D38*C38)"
$ws.Range("G38").Formula = "=(F38*1.21)"
$ws.Range("H38").Value = "https://arduino-shop.cz/arduino/1421-eses-i2c-20x4-display-pro-jednodeskove-pocitace.html?gclid=Cj0KCQjwuL_8BRCXARIsAGiC51AaUN-iQm0k_3qSD826rhlR5hUuxtvJMt2UjPC1cdwT4N_Vt1o1S2saAtzAEALw_wcB"

# --- Row 39 ---
$ws.Range("E39").Formula = "=(D39*1.21)"
$ws.Range("F39").Formula = "=(D39*C39)"
$ws.Range("G39").Formula = "=(F39*1.21)"

# --- Update the view: frozen pane top-left cell and active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("D39").Select()
